$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "rafatpc"
$ws.Range("C5").Value = "Liliyan Krumov"

$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/rafatpc") | Out-Null
$ws.Range("C7").Value = "https://github.com/rafatpc"

$ws.Range("C5:E5").Select()
